# Generate Report for Handback
# Adds a second handback record (afb46ee5-...) alongside the existing
# one (cc54d45f-... renamed to ac2f5080-...) across the Overview, zh-cn
# and de-de sheets, expanding each sheet's table by one row.

$wb = $excel.ActiveWorkbook

# ---- new identifiers used throughout ----
$oldGuid = "cc54d45f-e499-49a1-8ed4-e1bdfd642a56"
$guid1   = "ac2f5080-7f2a-49cf-be49-8ed8dacc307e"
$guid2   = "afb46ee5-6896-4257-a56d-04be8f8c5f92"

$zhHash1 = "6271b15e4e1671a6ee414920087270d9c3b9af42"
$deHash1 = "6271b15e4e1671a6ee414920087270d9c3b9af42"
$zhHash2 = "03a0aa0ab94907ca03c25347035788a6465ddd12"
$deHash2 = "03a0aa0ab94907ca03c25347035788a6465ddd12"

$dateOverview1 = "2016-08-23 00:57:52"
$dateZhHO1     = "2016-08-23 00:57:47"
$dateZhHB1     = "2016-08-23 00:58:09"
$dateDeHB1     = "2016-08-23 00:58:16"

# =====================================================================
# Sheet "Overview"
# =====================================================================
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)

# --- update row 2 (renamed guid, refreshed timestamp) ---
$ws.Range("A2").Value = ($guid1 + ".md")
$ws.Range("B2").Value = ("e2e\" + $guid1 + ".md")
$ws.Range("G2").Value = $dateOverview1

# --- add row 3 for the new handback file ---
$newRow = $lo.ListRows.Add()
$ws.Range("A3").Value = ($guid2 + ".md")
$ws.Range("B3").Value = ("e2e\" + $guid2 + ".md")
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = $dateOverview1

# visual formatting to mirror the hyperlink-styled B2 cell
$ws.Range("B3").Font.Underline = 1
$ws.Range("B3").Font.Color = 15570276

# --- hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("B3"), ("https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6fb85384755ea085329c38292c5454ed05d5e955/e2e/" + $guid2 + ".md"), "", "", ("e2e\" + $guid2 + ".md")) | Out-Null

# =====================================================================
# Sheet "zh-cn"
# =====================================================================
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)

# --- update row 2 ---
$ws.Range("A2").Value = ($guid1 + ".md")
$ws.Range("G2").Value = ($guid1 + "." + $zhHash1 + ".zh-cn.xlf")
$ws.Range("H2").Value = $dateZhHO1
$ws.Range("I2").Value = ($guid1 + ".md")
$ws.Range("J2").Value = ($guid1 + "." + $zhHash1 + ".zh-cn.xlf")
$ws.Range("K2").Value = $dateZhHB1

# --- add row 3 ---
$newRow = $lo.ListRows.Add()
$ws.Range("A3").Value = ($guid2 + ".md")
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = ($guid2 + "." + $zhHash2 + ".zh-cn.xlf")
$ws.Range("H3").Value = $dateZhHO1
$ws.Range("I3").Value = ($guid2 + ".md")
$ws.Range("J3").Value = ($guid2 + "." + $zhHash2 + ".zh-cn.xlf")
$ws.Range("K3").Value = $dateZhHB1
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

$ws.Range("A3").Font.Underline = 1
$ws.Range("A3").Font.Color = 15570276
$ws.Range("I3").Font.Underline = 1
$ws.Range("I3").Font.Color = 15570276
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Add($ws.Range("A3"), ("https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6fb85384755ea085329c38292c5454ed05d5e955/e2e/" + $guid2 + ".md"), "", "", ($guid2 + ".md")) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), ("https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d87301d650263f8551073e56557c6f0a1fc35084/e2e/" + $guid2 + ".md"), "", "", ($guid2 + ".md")) | Out-Null

# =====================================================================
# Sheet "de-de"
# =====================================================================
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)

# --- update row 2 ---
$ws.Range("A2").Value = ($guid1 + ".md")
$ws.Range("G2").Value = ($guid1 + "." + $deHash1 + ".de-de.xlf")
$ws.Range("I2").Value = ($guid1 + ".md")
$ws.Range("J2").Value = ($guid1 + "." + $deHash1 + ".de-de.xlf")
$ws.Range("K2").Value = $dateDeHB1

# --- add row 3 ---
$newRow = $lo.ListRows.Add()
$ws.Range("A3").Value = ($guid2 + ".md")
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = ($guid2 + "." + $deHash2 + ".de-de.xlf")
$ws.Range("H3").Value = "2016-08-23 00:55:33"
$ws.Range("I3").Value = ($guid2 + ".md")
$ws.Range("J3").Value = ($guid2 + "." + $deHash2 + ".de-de.xlf")
$ws.Range("K3").Value = $dateDeHB1
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

$ws.Range("A3").Font.Underline = 1
$ws.Range("A3").Font.Color = 15570276
$ws.Range("I3").Font.Underline = 1
$ws.Range("I3").Font.Color = 15570276
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Add($ws.Range("A3"), ("https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6fb85384755ea085329c38292c5454ed05d5e955/e2e/" + $guid2 + ".md"), "", "", ($guid2 + ".md")) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), ("https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/656959c2ba811bcbccc5d4d9b797ae135d95e994/e2e/" + $guid2 + ".md"), "", "", ($guid2 + ".md")) | Out-Null

Write-Output "done"
